# Update header K1 from "AGENT " to "VEH# TRAILER"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "VEH# TRAILER"

# Add trailer vehicle numbers for rows 5 and 9
$ws.Range("K5").Value = "TR20282"
$ws.Range("K9").Value = "TR272625"

# Update the active selection to match the recorded cursor position
$ws.Range("K18").Select()
